$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap mismatched match rows (F:V only; A and E are row-position metadata that stay put) ---
# Row 19
$ws.Range("F19").Value = "Talleres Cordoba"
$ws.Range("G19").Value = 2
$ws.Range("H19").Value = "Huracan"
$ws.Range("I19").Value = 1
$ws.Range("J19").Value = 1.74
$ws.Range("K19").Value = "21/08/2023 23:12"
$ws.Range("L19").Value = 1.94
$ws.Range("M19").Value = "26/08/2023 23:50"
$ws.Range("N19").Value = 3.44
$ws.Range("O19").Value = "21/08/2023 23:12"
$ws.Range("P19").Value = 3.23
$ws.Range("Q19").Value = "26/08/2023 23:50"
$ws.Range("R19").Value = 5.57
$ws.Range("S19").Value = "21/08/2023 23:12"
$ws.Range("T19").Value = 4.73
$ws.Range("U19").Value = "26/08/2023 23:50"
$ws.Range("V19").Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/talleres-cordoba-huracan/l4ycuMBq/"

# Row 20
$ws.Range("F20").Value = "Newells Old Boys"
$ws.Range("G20").Value = 1
$ws.Range("H20").Value = "Lanus"
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 2.19
$ws.Range("K20").Value = "22/08/2023 01:42"
$ws.Range("L20").Value = 2.09
$ws.Range("M20").Value = "26/08/2023 23:43"
$ws.Range("N20").Value = 2.94
$ws.Range("O20").Value = "22/08/2023 01:42"
$ws.Range("P20").Value = 3.01
$ws.Range("Q20").Value = "26/08/2023 23:43"
$ws.Range("R20").Value = 3.79
$ws.Range("S20").Value = "22/08/2023 01:42"
$ws.Range("T20").Value = 4.48
$ws.Range("U20").Value = "26/08/2023 23:42"
$ws.Range("V20").Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/newells-old-boys-lanus/OdYgrAqL/"

# Row 22
$ws.Range("F22").Value = "Platense"
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = "Defensa y Justicia"
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 2.61
$ws.Range("K22").Value = "21/08/2023 08:12"
$ws.Range("L22").Value = 2.08
$ws.Range("M22").Value = "27/08/2023 19:27"
$ws.Range("N22").Value = 3.05
$ws.Range("O22").Value = "21/08/2023 08:12"
$ws.Range("P22").Value = 3.11
$ws.Range("Q22").Value = "27/08/2023 19:22"
$ws.Range("R22").Value = 2.86
$ws.Range("S22").Value = "21/08/2023 08:12"
$ws.Range("T22").Value = 4.32
$ws.Range("U22").Value = "27/08/2023 19:27"
$ws.Range("V22").Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/platense-defensa-y-justicia/8U9Qxl7r/"

# Row 23
$ws.Range("F23").Value = "Godoy Cruz"
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = "Central Cordoba"
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 1.87
$ws.Range("K23").Value = "22/08/2023 01:42"
$ws.Range("L23").Value = 1.85
$ws.Range("M23").Value = "27/08/2023 19:26"
$ws.Range("N23").Value = 3.32
$ws.Range("O23").Value = "22/08/2023 01:42"
$ws.Range("P23").Value = 3.44
$ws.Range("Q23").Value = "27/08/2023 19:26"
$ws.Range("R23").Value = 4.39
$ws.Range("S23").Value = "22/08/2023 01:42"
$ws.Range("T23").Value = 4.87
$ws.Range("U23").Value = "27/08/2023 19:26"
$ws.Range("V23").Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/godoy-cruz-central-cordoba-santiago-del-estero/6LNbsUbR/"

# Row 60
$ws.Range("F60").Value = "Huracan"
$ws.Range("G60").Value = 2
$ws.Range("H60").Value = "Gimnasia L.P."
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 1.8
$ws.Range("K60").Value = "15/09/2023 01:13"
$ws.Range("L60").Value = 1.9
$ws.Range("M60").Value = "20/09/2023 01:55"
$ws.Range("N60").Value = 3.3
$ws.Range("O60").Value = "15/09/2023 01:13"
$ws.Range("P60").Value = 3.1
$ws.Range("Q60").Value = "20/09/2023 01:55"
$ws.Range("R60").Value = 4.82
$ws.Range("S60").Value = "15/09/2023 01:13"
$ws.Range("T60").Value = 5.38
$ws.Range("U60").Value = "20/09/2023 01:55"
$ws.Range("V60").Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/huracan-gimnasia-l-p/IZ3aGfe4/"

# Row 61
$ws.Range("F61").Value = "Rosario Central"
$ws.Range("G61").Value = 1
$ws.Range("H61").Value = "Independiente"
$ws.Range("I61").Value = 1
$ws.Range("J61").Value = 2.1
$ws.Range("K61").Value = "16/09/2023 01:13"
$ws.Range("L61").Value = 2.78
$ws.Range("M61").Value = "20/09/2023 01:55"
$ws.Range("N61").Value = 3.12
$ws.Range("O61").Value = "16/09/2023 01:13"
$ws.Range("P61").Value = 3.02
$ws.Range("Q61").Value = "20/09/2023 01:55"
$ws.Range("R61").Value = 4.04
$ws.Range("S61").Value = "16/09/2023 01:13"
$ws.Range("T61").Value = 2.93
$ws.Range("U61").Value = "20/09/2023 01:55"
$ws.Range("V61").Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/rosario-central-independiente/vq5eHztb/"

# Row 62
$ws.Range("F62").Value = "Defensa y Justicia"
$ws.Range("G62").Value = 2
$ws.Range("H62").Value = "Tigre"
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 2.49
$ws.Range("K62").Value = "15/09/2023 23:13"
$ws.Range("L62").Value = 2.67
$ws.Range("M62").Value = "20/09/2023 20:58"
$ws.Range("N62").Value = 3.2
$ws.Range("O62").Value = "15/09/2023 23:13"
$ws.Range("P62").Value = 3.03
$ws.Range("Q62").Value = "20/09/2023 20:53"
$ws.Range("R62").Value = 3.05
$ws.Range("S62").Value = "15/09/2023 23:13"
$ws.Range("T62").Value = 3.04
$ws.Range("U62").Value = "20/09/2023 20:53"
$ws.Range("V62").Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/defensa-y-justicia-tigre/0bg4FEAA/"

# Row 63
$ws.Range("F63").Value = "Barracas Central"
$ws.Range("G63").Value = 1
$ws.Range("H63").Value = "Banfield"
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 2.72
$ws.Range("K63").Value = "16/09/2023 23:12"
$ws.Range("L63").Value = 3.41
$ws.Range("M63").Value = "20/09/2023 20:59"
$ws.Range("N63").Value = 2.92
$ws.Range("O63").Value = "16/09/2023 23:12"
$ws.Range("P63").Value = 2.8
$ws.Range("Q63").Value = "20/09/2023 20:52"
$ws.Range("R63").Value = 3
$ws.Range("S63").Value = "16/09/2023 23:12"
$ws.Range("T63").Value = 2.61
$ws.Range("U63").Value = "20/09/2023 20:59"
$ws.Range("V63").Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/barracas-central-banfield/UDN4fDYp/"

# --- Append new rows 80-85 ---
# Copy formatting from an existing data row (row 2) to establish correct per-cell styles (A: s=1, E: s=2)
$ws.Range("A2:V2").Copy()
$ws.Range("A80:V80").PasteSpecial(-4122)
$ws.Range("A2:V2").Copy()
$ws.Range("A81:V81").PasteSpecial(-4122)
$ws.Range("A2:V2").Copy()
$ws.Range("A82:V82").PasteSpecial(-4122)
$ws.Range("A2:V2").Copy()
$ws.Range("A83:V83").PasteSpecial(-4122)
$ws.Range("A2:V2").Copy()
$ws.Range("A84:V84").PasteSpecial(-4122)
$ws.Range("A2:V2").Copy()
$ws.Range("A85:V85").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 80
$ws.Range("A80").Value = 79
$ws.Range("B80").Value = "argentina"
$ws.Range("C80").Value = "copa-de-la-liga-profesional"
$ws.Range("D80").Value = "'2023"
$ws.Range("D80").Style = "Normal"
$ws.Range("E80").Value = 45194.875
$ws.Range("F80").Value = "Platense"
$ws.Range("G80").Value = 1
$ws.Range("H80").Value = "Union de Santa Fe"
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 2.64
$ws.Range("K80").Value = "21/09/2023 22:42"
$ws.Range("L80").Value = 2.59
$ws.Range("M80").Value = "25/09/2023 20:59"
$ws.Range("N80").Value = 2.97
$ws.Range("O80").Value = "21/09/2023 22:42"
$ws.Range("P80").Value = 2.77
$ws.Range("Q80").Value = "25/09/2023 20:58"
$ws.Range("R80").Value = 2.9
$ws.Range("S80").Value = "21/09/2023 22:42"
$ws.Range("T80").Value = 3.5
$ws.Range("U80").Value = "25/09/2023 20:59"
$ws.Range("V80").Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/platense-union-de-santa-fe/nFwHsXg3/"

# Row 81
$ws.Range("A81").Value = 80
$ws.Range("B81").Value = "argentina"
$ws.Range("C81").Value = "copa-de-la-liga-profesional"
$ws.Range("D81").Value = "'2023"
$ws.Range("D81").Style = "Normal"
$ws.Range("E81").Value = 45194.875
$ws.Range("F81").Value = "Sarmiento Junin"
$ws.Range("G81").Value = 0
$ws.Range("H81").Value = "Belgrano"
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 2.37
$ws.Range("K81").Value = "21/09/2023 21:12"
$ws.Range("L81").Value = 3.16
$ws.Range("M81").Value = "25/09/2023 20:58"
$ws.Range("N81").Value = 2.97
$ws.Range("O81").Value = "21/09/2023 21:12"
$ws.Range("P81").Value = 2.72
$ws.Range("Q81").Value = "25/09/2023 20:58"
$ws.Range("R81").Value = 3.51
$ws.Range("S81").Value = "21/09/2023 21:12"
$ws.Range("T81").Value = 2.86
$ws.Range("U81").Value = "25/09/2023 20:58"
$ws.Range("V81").Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/sarmiento-junin-ca-belgrano-de-cordoba/txwLti99/"

# Row 82
$ws.Range("A82").Value = 81
$ws.Range("B82").Value = "argentina"
$ws.Range("C82").Value = "copa-de-la-liga-profesional"
$ws.Range("D82").Value = "'2023"
$ws.Range("D82").Style = "Normal"
$ws.Range("E82").Value = 45194.97916666666
$ws.Range("F82").Value = "Colon Santa Fe"
$ws.Range("G82").Value = 3
$ws.Range("H82").Value = "Argentinos Jrs"
$ws.Range("I82").Value = 1
$ws.Range("J82").Value = 3.03
$ws.Range("K82").Value = "21/09/2023 23:42"
$ws.Range("L82").Value = 2.95
$ws.Range("M82").Value = "25/09/2023 23:25"
$ws.Range("N82").Value = 3.09
$ws.Range("O82").Value = "21/09/2023 23:42"
$ws.Range("P82").Value = 3.16
$ws.Range("Q82").Value = "25/09/2023 23:20"
$ws.Range("R82").Value = 2.46
$ws.Range("S82").Value = "21/09/2023 23:42"
$ws.Range("T82").Value = 2.64
$ws.Range("U82").Value = "25/09/2023 23:29"
$ws.Range("V82").Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/colon-santa-fe-argentinos-jrs/0IHoSifG/"

# Row 83
$ws.Range("A83").Value = 82
$ws.Range("B83").Value = "argentina"
$ws.Range("C83").Value = "copa-de-la-liga-profesional"
$ws.Range("D83").Value = "'2023"
$ws.Range("D83").Style = "Normal"
$ws.Range("E83").Value = 45194.97916666666
$ws.Range("F83").Value = "Godoy Cruz"
$ws.Range("G83").Value = 1
$ws.Range("H83").Value = "Racing Club"
$ws.Range("I83").Value = 1
$ws.Range("J83").Value = 2.34
$ws.Range("K83").Value = "24/09/2023 05:42"
$ws.Range("L83").Value = 2.29
$ws.Range("M83").Value = "25/09/2023 23:29"
$ws.Range("N83").Value = 3.2
$ws.Range("O83").Value = "24/09/2023 05:42"
$ws.Range("P83").Value = 3.36
$ws.Range("Q83").Value = "25/09/2023 23:29"
$ws.Range("R83").Value = 3.3
$ws.Range("S83").Value = "24/09/2023 05:42"
$ws.Range("T83").Value = 3.35
$ws.Range("U83").Value = "25/09/2023 23:29"
$ws.Range("V83").Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/godoy-cruz-racing-club/8OvDrDvc/"

# Row 84
$ws.Range("A84").Value = 83
$ws.Range("B84").Value = "argentina"
$ws.Range("C84").Value = "copa-de-la-liga-profesional"
$ws.Range("D84").Value = "'2023"
$ws.Range("D84").Style = "Normal"
$ws.Range("E84").Value = 45195.08333333334
$ws.Range("F84").Value = "Talleres Cordoba"
$ws.Range("G84").Value = 4
$ws.Range("H84").Value = "Barracas Central"
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 1.52
$ws.Range("K84").Value = "21/09/2023 23:42"
$ws.Range("L84").Value = 1.67
$ws.Range("M84").Value = "26/09/2023 01:58"
$ws.Range("N84").Value = 4.1
$ws.Range("O84").Value = "21/09/2023 23:42"
$ws.Range("P84").Value = 3.69
$ws.Range("Q84").Value = "26/09/2023 01:58"
$ws.Range("R84").Value = 6.9
$ws.Range("S84").Value = "21/09/2023 23:42"
$ws.Range("T84").Value = 5.96
$ws.Range("U84").Value = "26/09/2023 01:58"
$ws.Range("V84").Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/talleres-cordoba-barracas-central/v9GkRB9M/"

# Row 85
$ws.Range("A85").Value = 84
$ws.Range("B85").Value = "argentina"
$ws.Range("C85").Value = "copa-de-la-liga-profesional"
$ws.Range("D85").Value = "'2023"
$ws.Range("D85").Style = "Normal"
$ws.Range("E85").Value = 45195.08333333334
$ws.Range("F85").Value = "Atl. Tucuman"
$ws.Range("G85").Value = 1
$ws.Range("H85").Value = "Arsenal Sarandi"
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 1.74
$ws.Range("K85").Value = "22/09/2023 01:12"
$ws.Range("L85").Value = 1.72
$ws.Range("M85").Value = "26/09/2023 01:58"
$ws.Range("N85").Value = 3.47
$ws.Range("O85").Value = "22/09/2023 01:12"
$ws.Range("P85").Value = 3.54
$ws.Range("Q85").Value = "26/09/2023 01:58"
$ws.Range("R85").Value = 5.5
$ws.Range("S85").Value = "22/09/2023 01:12"
$ws.Range("T85").Value = 5.8
$ws.Range("U85").Value = "26/09/2023 01:58"
$ws.Range("V85").Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/atl-tucuman-arsenal-sarandi/bgcMLkXq/"
